$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text cells (coin names & links) ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

# --- Numeric-looking text cells (price/volume) must stay as Text, not auto-convert to Number/Percent ---
$numericCells = @{
    "D2" = "245.38"
    "D3" = "28.40"
    "E3" = "-3.09%"
    "D4" = "5.289"
    "E4" = "1.98%"
    "E5" = "-0.41%"
    "D6" = "6.642"
    "E6" = "1.49%"
    "D7" = "3.207"
    "E7" = "3.75%"
    "D8" = "0.8501"
    "E8" = "-1.03%"
    "D9" = "0.8824"
    "E9" = "1.41%"
    "D10" = "0.1381"
    "E10" = "1.21%"
    "D11" = "0.07082"
    "E11" = "0.10%"
    "D12" = "0.03144"
    "E12" = "5.11%"
    "D13" = "0.09221"
    "E13" = "-1.73%"
    "D14" = "0.001526"
    "E14" = "0.06%"
    "D15" = "0.0005957"
    "E15" = "-94.20%"
    "D16" = "0.006041"
    "E16" = "1.05%"
    "D17" = "3.494"
    "E17" = "0.09%"
    "E18" = "0.21%"
    "E19" = "-0.54%"
    "D20" = "0.03315"
    "E20" = "-1.19%"
    "D21" = "0.1294"
    "E21" = "-0.21%"
    "D22" = "3.521"
    "E22" = "1.02%"
    "D23" = "0.04072"
    "E23" = "-1.45%"
    "E24" = "-0.07%"
    "D25" = "0.001219"
    "E25" = "-0.55%"
    "D26" = "0.004155"
    "E26" = "-16.96%"
    "E27" = "-0.83%"
    "D28" = "0.0001448"
    "D40" = "0.03782"
    "E40" = "0.67%"
    "D41" = "0.1067"
    "E41" = "-0.64%"
    "D42" = "0.003740"
    "E42" = "-35.17%"
    "D43" = "0.002239"
    "E43" = "-7.78%"
    "D44" = "0.009473"
    "E44" = "0.34%"
    "D45" = "0.00005267"
    "E45" = "0.08%"
    "E46" = "0.00%"
    "D47" = "0.08906"
    "E47" = "56.31%"
    "D48" = "0.002269"
    "E48" = "0.42%"
    "D49" = "0.00002099"
    "E49" = "0.00%"
    "D50" = "0.0001999"
    "E50" = "0.00%"
}

foreach ($ref in $numericCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericCells[$ref]
    $cell.Style = "Normal"
}
